$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be pre-formatted as Text
# so Excel stores them as strings (matching the source inlineStr cells), then the
# style is reset back to Normal so no stray number-format style is left behind.
$textCells = @("D5", "D6", "D9", "D11", "D14", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D30", "D32", "D33", "D36", "D37", "D40", "D41", "D42", "D43", "D45", "D49", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.514.13"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "3.175.90"
$ws.Range("E3").Value = "  -3.85%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "586.66"
$ws.Range("E5").Value = "  -2.90%  "
$ws.Range("D6").Value = "135.08"
$ws.Range("E6").Value = "  -5.16%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.172.08"
$ws.Range("E8").Value = "  -3.94%  "
$ws.Range("D9").Value = "0.507"
$ws.Range("E9").Value = "  -2.49%  "
$ws.Range("E10").Value = "  -5.55%  "
$ws.Range("D11").Value = "5.25"
$ws.Range("E11").Value = "  -4.02%  "
$ws.Range("E12").Value = "  -3.81%  "
$ws.Range("E13").Value = "  -5.15%  "
$ws.Range("D14").Value = "33.24"
$ws.Range("E14").Value = "  -4.11%  "
$ws.Range("D15").Value = "3.702.86"
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").Value = "3.182.51"
$ws.Range("E17").Value = "  -3.61%  "
$ws.Range("D18").Value = "62.509.36"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").Value = "6.50"
$ws.Range("E19").Value = "  -5.31%  "
$ws.Range("D20").Value = "455.12"
$ws.Range("E20").Value = "  -5.35%  "
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").Value = "0.700"
$ws.Range("E22").Value = "  -4.21%  "
$ws.Range("D23").Value = "7.59"
$ws.Range("E23").Value = "  -5.42%  "
$ws.Range("D24").Value = "83.59"
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("D25").Value = "13.21"
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -6.63%  "
$ws.Range("D30").Value = "7.71"
$ws.Range("E30").Value = "  -5.32%  "
$ws.Range("E31").Value = "  -7.41%  "
$ws.Range("D32").Value = "27.28"
$ws.Range("E32").Value = "  -6.44%  "
$ws.Range("D33").Value = "0.104"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("E34").Value = "  -6.68%  "
$ws.Range("E35").Value = "  -6.15%  "
$ws.Range("D36").Value = "5.90"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").Value = "51.04"
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("D38").Value = "0.0₃0698"
$ws.Range("E38").Value = "  -6.73%  "
$ws.Range("D40").Value = "2.72"
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.112"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "7.97"
$ws.Range("E42").Value = "  -4.85%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "393.25"
$ws.Range("E43").Value = "  -7.46%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.805.20"
$ws.Range("E44").Value = "  -8.32%  "
$ws.Range("D45").Value = "36.24"
$ws.Range("E45").Value = "  +3.94%  "
$ws.Range("E46").Value = "  -6.34%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("D49").Value = "124.77"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "25.22"
$ws.Range("E50").Value = "  -4.07%  "
$ws.Range("E51").Value = "  -3.87%  "

foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
